# Nov 4 Team Meeting Attendance
# Add a new "Nov 4" meeting column (K) on the TEAM sheet and mark
# attendance for everyone who was present.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TEAM")

$check = [char]0x2714

# Header: month + day for the new Nov 4 meeting (mirrors columns C:J)
$ws.Range("K3").Value = "Nov"
$ws.Range("K4").Value = 4

# Attendance checkmarks for everyone who attended the Nov 4 meeting
$ws.Range("K5").Value = $check
$ws.Range("K6").Value = $check
$ws.Range("K7").Value = $check
$ws.Range("K8").Value = $check
$ws.Range("K9").Value = $check
$ws.Range("K10").Value = $check
$ws.Range("K11").Value = $check

# Update selections left over from editing: collapse the multi-range
# selections on SPONSOR/TA back to a single cell, and land on T6 on TEAM.
$sponsor = $wb.Worksheets.Item("SPONSOR")
$null = $sponsor.Range("F11").Select()

$ta = $wb.Worksheets.Item("TA")
$null = $ta.Range("F5").Select()

$null = $ws.Activate()
$null = $ws.Range("T6").Select()
